$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024745733218082
$ws.Range("D2").Value = 1.049606799248559
$ws.Range("E2").Value = 1.036776901417529
$ws.Range("F2").Value = 1.051638280568072
$ws.Range("I2").Value = 1.039231582915969
$ws.Range("J2").Value = 1.029918625759161
$ws.Range("K2").Value = 1.052363388078976
$ws.Range("L2").Value = 1.039569674389972
$ws.Range("M2").Value = 1.054389227793618
$ws.Range("N2").Value = 1.014014349005517
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025679754426089
$ws.Range("D3").Value = 1.050292319505274
$ws.Range("E3").Value = 1.037593648790926
$ws.Range("F3").Value = 1.052510448914736
$ws.Range("I3").Value = 1.039382792865286
$ws.Range("J3").Value = 1.030491765034563
$ws.Range("K3").Value = 1.052861515574064
$ws.Range("L3").Value = 1.040196074380964
$ws.Range("M3").Value = 1.055073929407504
$ws.Range("N3").Value = 1.01420451217778
$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.026284387139818
$ws.Range("D4").Value = 1.050733584677209
$ws.Range("E4").Value = 1.038122176026911
$ws.Range("F4").Value = 1.053073575717682
$ws.Range("I4").Value = 1.039477811958181
$ws.Range("J4").Value = 1.030862283780101
$ws.Range("K4").Value = 1.053180851901997
$ws.Range("L4").Value = 1.040600785779166
$ws.Range("M4").Value = 1.055515116286644
$ws.Range("N4").Value = 1.014327415501061
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026538635791623
$ws.Range("D5").Value = 1.050918536960051
$ws.Range("E5").Value = 1.038344376388188
$ws.Range("F5").Value = 1.053310019406072
$ws.Range("I5").Value = 1.039517080790841
$ws.Range("J5").Value = 1.03101796743739
$ws.Range("K5").Value = 1.053314384312616
$ws.Range("L5").Value = 1.040770778979461
$ws.Range("M5").Value = 1.055700144144649
$ws.Range("N5").Value = 1.014379049041113
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026581328815064
$ws.Range("D6").Value = 1.050949558621761
$ws.Range("E6").Value = 1.038381685206994
$ws.Range("F6").Value = 1.05334970200362
$ws.Range("I6").Value = 1.039523634455435
$ws.Range("J6").Value = 1.031044102549844
$ws.Range("K6").Value = 1.053336762916419
$ws.Range("L6").Value = 1.040799312904771
$ws.Range("M6").Value = 1.055731184869974
$ws.Range("N6").Value = 1.01438771647711
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026287784183902
$ws.Range("D7").Value = 1.050736058204455
$ws.Range("E7").Value = 1.038125145051266
$ws.Range("F7").Value = 1.053076736250477
$ws.Range("I7").Value = 1.039478339332684
$ws.Range("J7").Value = 1.030864364357733
$ws.Range("K7").Value = 1.053182638986734
$ws.Range("L7").Value = 1.040603057817562
$ws.Range("M7").Value = 1.055517590396889
$ws.Range("N7").Value = 1.014328105568946
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025061336028103
$ws.Range("D8").Value = 1.049838951592644
$ws.Range("E8").Value = 1.037052916526958
$ws.Range("F8").Value = 1.051933286469394
$ws.Range("I8").Value = 1.039283269037104
$ws.Range("J8").Value = 1.030112390889133
$ws.Range("K8").Value = 1.052532349658661
$ws.Range("L8").Value = 1.039781495096846
$ws.Range("M8").Value = 1.054621010251458
$ws.Range("N8").Value = 1.014078645369558
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022902193796433
$ws.Range("D9").Value = 1.048240513146687
$ws.Range("E9").Value = 1.035163855848316
$ws.Range("F9").Value = 1.049909076536553
$ws.Range("I9").Value = 1.038917952498597
$ws.Range("J9").Value = 1.02878474549903
$ws.Range("K9").Value = 1.051363671099026
$ws.Range("L9").Value = 1.038329163650463
$ws.Range("M9").Value = 1.053026938859836
$ws.Range("N9").Value = 1.013637967888356
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021464167088691
$ws.Range("D10").Value = 1.047163152865996
$ws.Range("E10").Value = 1.033904789756183
$ws.Range("F10").Value = 1.048553440379353
$ws.Range("I10").Value = 1.038659964896189
$ws.Range("J10").Value = 1.027897967487849
$ws.Range("K10").Value = 1.050569335272481
$ws.Range("L10").Value = 1.037357891010396
$ws.Range("M10").Value = 1.051954778524046
$ws.Range("N10").Value = 1.013343463122549
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020841826904105
$ws.Range("D11").Value = 1.046693883290853
$ws.Range("E11").Value = 1.03335968871655
$ws.Range("F11").Value = 1.047964991684629
$ws.Range("I11").Value = 1.038544840510274
$ws.Range("J11").Value = 1.027513593657428
$ws.Range("K11").Value = 1.050221790361723
$ws.Range("L11").Value = 1.036936608022751
$ws.Range("M11").Value = 1.051488298490576
$ws.Range("N11").Value = 1.013215772106622
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020610712985719
$ws.Range("D12").Value = 1.046519161888694
$ws.Range("E12").Value = 1.033157227550214
$ws.Range("F12").Value = 1.047746199263504
$ws.Range("I12").Value = 1.038501566209243
$ws.Range("J12").Value = 1.027370761752021
$ws.Range("K12").Value = 1.050092158529936
$ws.Range("L12").Value = 1.036780018204436
$ws.Range("M12").Value = 1.051314693679729
$ws.Range("N12").Value = 1.013168316924544
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020660285390849
$ws.Range("D13").Value = 1.046556658920619
$ws.Range("E13").Value = 1.033200655522154
$ws.Range("F13").Value = 1.04779314075043
$ws.Range("I13").Value = 1.038510871846914
$ws.Range("J13").Value = 1.027401402309829
$ws.Range("K13").Value = 1.050119989327629
$ws.Range("L13").Value = 1.036813612059138
$ws.Range("M13").Value = 1.051351947565453
$ws.Range("N13").Value = 1.013178497351721
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020822721912856
$ws.Range("D14").Value = 1.046679449196768
$ws.Range("E14").Value = 1.033342952928587
$ws.Range("F14").Value = 1.047946910641071
$ws.Range("I14").Value = 1.038541273881887
$ws.Range("J14").Value = 1.027501788313854
$ws.Range("K14").Value = 1.050211085924103
$ws.Range("L14").Value = 1.036923666427923
$ws.Range("M14").Value = 1.051473955055704
$ws.Range("N14").Value = 1.01321184995604
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020922811183189
$ws.Range("D15").Value = 1.046755049557446
$ws.Range("E15").Value = 1.033430628898143
$ws.Range("F15").Value = 1.048041624714135
$ws.Range("I15").Value = 1.038559937764736
$ws.Range("J15").Value = 1.027563631733045
$ws.Range("K15").Value = 1.050267142274687
$ws.Range("L15").Value = 1.036991460464477
$ws.Range("M15").Value = 1.051549083762683
$ws.Range("N15").Value = 1.013232396284331
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021505476866041
$ws.Range("D16").Value = 1.047194238622486
$ws.Range("E16").Value = 1.033940968164234
$ws.Range("F16").Value = 1.048592463363565
$ws.Range("I16").Value = 1.038667533473614
$ws.Range("J16").Value = 1.027923468891651
$ws.Range("K16").Value = 1.05059232512948
$ws.Range("L16").Value = 1.037385835180243
$ws.Range("M16").Value = 1.051985690448207
$ws.Range("N16").Value = 1.01335193402701
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021871057753944
$ws.Range("D17").Value = 1.047468991135284
$ws.Range("E17").Value = 1.034261113531735
$ws.Range("F17").Value = 1.048937602853318
$ws.Range("I17").Value = 1.038734112055971
$ws.Range("J17").Value = 1.028149080494599
$ws.Range("K17").Value = 1.050795343079931
$ws.Range("L17").Value = 1.037633025123891
$ws.Range("M17").Value = 1.052258966846001
$ws.Range("N17").Value = 1.013426871967402
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022084327169375
$ws.Range("D18").Value = 1.047628982689371
$ws.Range("E18").Value = 1.034447856827681
$ws.Range("F18").Value = 1.049138776953104
$ws.Range("I18").Value = 1.038772616748315
$ws.Range("J18").Value = 1.028280637950342
$ws.Range("K18").Value = 1.050913413315588
$ws.Range("L18").Value = 1.037777137736096
$ws.Range("M18").Value = 1.052418149335159
$ws.Range("N18").Value = 1.013470565724096
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022157051982807
$ws.Range("D19").Value = 1.047683490345285
$ws.Range("E19").Value = 1.034511532772886
$ws.Range("F19").Value = 1.049207348346532
$ws.Range("I19").Value = 1.038785689940682
$ws.Range("J19").Value = 1.028325489180012
$ws.Range("K19").Value = 1.050953613361568
$ws.Range("L19").Value = 1.037826264663236
$ws.Range("M19").Value = 1.05247238994063
$ws.Range("N19").Value = 1.013485461397204
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021831831042157
$ws.Range("D20").Value = 1.047439540389094
$ws.Range("E20").Value = 1.03422676413584
$ws.Range("F20").Value = 1.048900587115484
$ws.Range("I20").Value = 1.038727002875667
$ws.Range("J20").Value = 1.028124878427511
$ws.Range("K20").Value = 1.050773597024798
$ws.Range("L20").Value = 1.037606511122717
$ws.Range("M20").Value = 1.052229669104652
$ws.Range("N20").Value = 1.013418833517182
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020774886994692
$ws.Range("D21").Value = 1.046643301930103
$ws.Range("E21").Value = 1.033301049503434
$ws.Range("F21").Value = 1.047901635203118
$ws.Range("I21").Value = 1.038532335365988
$ws.Range("J21").Value = 1.027472228742857
$ws.Range("K21").Value = 1.050184275091314
$ws.Range("L21").Value = 1.03689126109344
$ws.Range("M21").Value = 1.051438036082738
$ws.Range("N21").Value = 1.01320202913486
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020110639324159
$ws.Range("D22").Value = 1.046140281966011
$ws.Range("E22").Value = 1.032719095402777
$ws.Range("F22").Value = 1.047272303257636
$ws.Range("I22").Value = 1.038406978476819
$ws.Range("J22").Value = 1.027061544725953
$ws.Range("K22").Value = 1.049810631993413
$ws.Range("L22").Value = 1.036440938945954
$ws.Range("M22").Value = 1.050938376905903
$ws.Range("N22").Value = 1.013065570791359
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020462741413164
$ws.Range("D23").Value = 1.046407168664822
$ws.Range("E23").Value = 1.033027592365352
$ws.Range("F23").Value = 1.047606042184994
$ws.Range("I23").Value = 1.038473712962144
$ws.Range("J23").Value = 1.027279287828542
$ws.Range("K23").Value = 1.050009001850698
$ws.Range("L23").Value = 1.036679721323571
$ws.Range("M23").Value = 1.051203438022099
$ws.Range("N23").Value = 1.01313792358189
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02184955580497
$ws.Range("D24").Value = 1.047452848738414
$ws.Range("E24").Value = 1.034242285124858
$ws.Range("F24").Value = 1.048917313366415
$ws.Range("I24").Value = 1.03873021622624
$ws.Range("J24").Value = 1.028135814417198
$ws.Range("K24").Value = 1.050783424202306
$ws.Range("L24").Value = 1.037618491872079
$ws.Range("M24").Value = 1.05224290815719
$ws.Range("N24").Value = 1.01342246579712
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023460140124081
$ws.Range("D25").Value = 1.048655824234922
$ws.Range("E25").Value = 1.035652174647054
$ws.Range("F25").Value = 1.050433476276608
$ws.Range("I25").Value = 1.039014945910028
$ws.Range("J25").Value = 1.02912827343809
$ws.Range("K25").Value = 1.051668491952341
$ws.Range("L25").Value = 1.038705168576207
$ws.Range("M25").Value = 1.053440715872997
$ws.Range("N25").Value = 1.013752021718461
